$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 & 3: shift start/stop times by 12 hours (AM -> PM) ---
$ws.Range("B2").Value = 0.7729166666666667
$ws.Range("C2").Value = 0.78055555555555556
$ws.Range("B3").Value = 0.58333333333333337
$ws.Range("C3").Value = 0.60763888888888895

# --- Row 5 content first (chronologically this was entered first: enemy ghost work) ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A5:C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A5").Value = "2/16/2020"
$ws.Range("B5").Value = 0.45833333333333331
$ws.Range("C5").Value = 0.70833333333333337
$ws.Range("D5").Value = 0.0625
$ws.Range("E5").Value = "4 hours 30 minutes"
$ws.Range("F5").Value = "Research/Programming"
$ws.Range("G5").Value = "Set up enemy class with enemy ghost that follows player"

# --- Row 6 content next ---
$ws.Range("A5:D5").Copy() | Out-Null
$ws.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A6").Value = "2/23/2020"
$ws.Range("B6").Value = 0.5
$ws.Range("C6").Value = 0.75
$ws.Range("D6").Value = 0.0625
$ws.Range("E6").Value = "4 hours 30 minutes"
$ws.Range("F6").Value = "Research/Programming"
$ws.Range("G6").Value = "Created sprite for protagonist and created animation classes"

# --- Row 4 content last (timesheet entry for earlier date, filled in afterward) ---
$ws.Range("A5:D5").Copy() | Out-Null
$ws.Range("A4:D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "2/7/2020"
$ws.Range("B4").Value = 0.75
$ws.Range("C4").Value = 0.91666666666666663
$ws.Range("D4").Value = 0.0625
$ws.Range("E4").Value = "2 hours 30 minutes"
$ws.Range("F4").Value = "Research/Setup"
$ws.Range("G4").Value = "Tried and failed to set up monogame extended third party library"

# --- Column widths (new col B, widened F & G to fit the new longer content) ---
$ws.Columns("B").ColumnWidth = 8.0013
$ws.Columns("F").ColumnWidth = 21.3346
$ws.Columns("G").ColumnWidth = 58.8346

# --- Selection ---
$ws.Range("C2").Select() | Out-Null
